# Applies the "Automatic update of files" edit to the Artfynd sheet:
#   1. Row 12 and row 13 swap their entire contents (observation order changed).
#   2. Four brand-new observation rows (19-22) are appended.
#   3. The sheet's used range grows from A1:AY18 to A1:AY22 (handled
#      automatically by Excel once the new cells are populated).
#
# NOTE: this COM engine only binds *positional* function parameters
# (named "-param value" binding silently fails), so every helper below
# is called positionally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a flat list of (ColumnLetter, Value) pairs into one row.
# Plain .Value assignment is used; Excel already keeps numbers/booleans/
# ordinary strings in their correct type for every column this sheet
# uses (verified: ints/floats stay numeric, "16:14"-style strings stay
# text, Booleans become t="b"). The only cells that need special
# handling are genuine date-look-alike strings and numeric-look-alike
# strings, which are routed through Set-ForcedText instead (see below).
# ---------------------------------------------------------------------
function Set-RowCells($ws, $rowNum, $pairs) {
    for ($i = 0; $i -lt $pairs.Length; $i += 2) {
        $col = $pairs[$i]
        $val = $pairs[$i + 1]
        $ws.Range("$col$rowNum").Value = $val
    }
}

# ---------------------------------------------------------------------
# Helper: force a value to be stored as literal text even though it
# looks like a date ("2023-09-26") or a number ("25"), which Excel
# would otherwise auto-convert when assigned straight to .Value.
# Done by staging the text (as Text-formatted) in a scratch cell, then
# copying only the *value* (PasteSpecial values) into the destination,
# so the destination cell keeps the default "General" style.
# ---------------------------------------------------------------------
function Set-ForcedText($scratch, $addr, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $scratch.Worksheet.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

$scratch = $ws.Range("A23")

# --- Row 12 (becomes what used to be row 13's record) -----------------
Set-RowCells $ws 12 @(
    "A", 112330621,
    "B", 77388,
    "C", "Ovaliderad",
    "D", "NT",
    "E", 6446,
    "F", "Kolflarnlav",
    "G", "Carbonicola anthracophila",
    "H", "(Nyl.) Bendiksby & Timdal",
    "I", "",
    "K", "",
    "P", "Bodhöjden (Bodhöjden), Ång",
    "Q", 577256,
    "R", 7056659,
    "S", 25,
    "T", "Västernorrland",
    "U", "Sollefteå",
    "V", "Ångermanland",
    "W", "Ramsele",
    "Z", "15:19",
    "AB", "15:19",
    "AD", $false,
    "AE", $false,
    "AG", $false,
    "AT", "",
    "AW", "Kamilla Andersson",
    "AX", "Kamilla Andersson",
    "AY", ""
)
Set-ForcedText $scratch "Y12" "2023-09-26"
Set-ForcedText $scratch "AA12" "2023-09-26"

# --- Row 13 (becomes what used to be row 12's record) -----------------
Set-RowCells $ws 13 @(
    "A", 112331359,
    "B", 77636,
    "C", "Ovaliderad",
    "D", "NT",
    "E", 6425,
    "F", "Garnlav",
    "G", "Alectoria sarmentosa",
    "H", "(Ach.) Ach.",
    "I", "",
    "K", "",
    "P", "Bodhöjden (Bodhöjden), Ång",
    "Q", 577102,
    "R", 7056694,
    "S", 25,
    "T", "Västernorrland",
    "U", "Sollefteå",
    "V", "Ångermanland",
    "W", "Ramsele",
    "Z", "16:14",
    "AB", "16:14",
    "AD", $false,
    "AE", $false,
    "AG", $false,
    "AT", "",
    "AW", "Kim Hultgren",
    "AX", "Kim Hultgren",
    "AY", ""
)
Set-ForcedText $scratch "Y13" "2023-09-26"
Set-ForcedText $scratch "AA13" "2023-09-26"

# --- Row 19 (new) -------------------------------------------------------
Set-RowCells $ws 19 @(
    "A", 112502528,
    "B", 56430,
    "C", "Ovaliderad",
    "D", "NT",
    "E", 100109,
    "F", "Tretåig hackspett",
    "G", "Picoides tridactylus",
    "H", "(Linnaeus, 1758)",
    "I", "",
    "K", "",
    "M", "gammalt bo",
    "P", "Bodhöjden (Bodhöjden), Ång",
    "Q", 577213,
    "R", 7056687,
    "S", 25,
    "T", "Västernorrland",
    "U", "Sollefteå",
    "V", "Ångermanland",
    "W", "Ramsele",
    "Z", "16:28",
    "AB", "16:28",
    "AD", $false,
    "AE", $false,
    "AG", $false,
    "AT", "",
    "AW", "Kamilla Andersson",
    "AX", "Kamilla Andersson",
    "AY", ""
)
Set-ForcedText $scratch "Y19" "2023-10-03"
Set-ForcedText $scratch "AA19" "2023-10-03"

# --- Row 20 (new) -------------------------------------------------------
Set-RowCells $ws 20 @(
    "A", 112502371,
    "B", 96735,
    "C", "Ovaliderad",
    "D", "VU",
    "E", 220787,
    "F", "Knärot",
    "G", "Goodyera repens",
    "H", "(L.) R. Br.",
    "K", "",
    "P", "Bodhöjden (Bodhöjden), Ång",
    "Q", 577167,
    "R", 7056735,
    "S", 25,
    "T", "Västernorrland",
    "U", "Sollefteå",
    "V", "Ångermanland",
    "W", "Ramsele",
    "Z", "16:17",
    "AB", "16:17",
    "AD", $false,
    "AE", $false,
    "AG", $false,
    "AT", "",
    "AW", "Kamilla Andersson",
    "AX", "Kamilla Andersson",
    "AY", ""
)
Set-ForcedText $scratch "Y20" "2023-10-03"
Set-ForcedText $scratch "AA20" "2023-10-03"
Set-ForcedText $scratch "I20" "25"

# --- Row 21 (new) -------------------------------------------------------
Set-RowCells $ws 21 @(
    "A", 112502199,
    "B", 56575,
    "C", "Ovaliderad",
    "D", "NT",
    "E", 103021,
    "F", "Talltita",
    "G", "Poecile montanus",
    "H", "(Conrad von Baldenstein, 1827)",
    "I", "",
    "K", "",
    "M", "födosökande",
    "P", "Bodhöjden (Bodhöjden), Ång",
    "Q", 577202,
    "R", 7056615,
    "S", 25,
    "T", "Västernorrland",
    "U", "Sollefteå",
    "V", "Ångermanland",
    "W", "Ramsele",
    "Z", "16:00",
    "AB", "16:00",
    "AD", $false,
    "AE", $false,
    "AG", $false,
    "AT", "",
    "AW", "Kamilla Andersson",
    "AX", "Kamilla Andersson",
    "AY", ""
)
Set-ForcedText $scratch "Y21" "2023-10-03"
Set-ForcedText $scratch "AA21" "2023-10-03"

# --- Row 22 (new) -------------------------------------------------------
Set-RowCells $ws 22 @(
    "A", 112502131,
    "B", 89834,
    "C", "Ovaliderad",
    "D", "NT",
    "E", 658,
    "F", "Rosenticka",
    "G", "Rhodofomes roseus",
    "H", "(Alb. & Schwein.) Kotl. & Pouzar",
    "I", "",
    "K", "",
    "P", "Sollefteå (Sollefteå), Ång",
    "Q", 577255,
    "R", 7056664,
    "S", 25,
    "T", "Västernorrland",
    "U", "Sollefteå",
    "V", "Ångermanland",
    "W", "Ramsele",
    "Z", "15:56",
    "AB", "15:56",
    "AD", $false,
    "AE", $false,
    "AG", $false,
    "AT", "",
    "AW", "Kim Hultgren",
    "AX", "Kim Hultgren",
    "AY", ""
)
Set-ForcedText $scratch "Y22" "2023-10-03"
Set-ForcedText $scratch "AA22" "2023-10-03"

# Drop the scratch row entirely so it never taints the sheet's used range.
$scratch.EntireRow.Delete()
